# Gate_Closure_Trigger.xlsx update:
# "Updating CPRA slide decks to include additional stations."
#
# 1) Station "17StCanal"/row 8 (WBV90, C8) measured gate closure value
#    updated from 0 to 2.5.
# 2) Three new monitoring stations appended to the bottom of the table
#    (rows 31-33): Venice, BayouSale, BayouBoeuf - each formatted the
#    same way as the other "MS River / GIWW" style rows above them
#    (row 23's format: column A uses the 0.00000000/Arial style, column
#    C uses the 0.0 style, column B is unstyled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing data value ---------------------------------
$ws.Cells.Item(8, 3).Value = 2.5

# --- 2) Append the three new station rows ---------------------------
# Clone the formatting of row 23 (an existing "new station" style row)
# down into rows 31:33 so the new cells pick up the same number formats
# / fonts/styles as their siblings instead of creating new style
# entries.
$ws.Range("A23").Copy()
$ws.Range("A31:A33").PasteSpecial(-4122)

$ws.Range("B23").Copy()
$ws.Range("B31:B33").PasteSpecial(-4122)

$ws.Range("C23").Copy()
$ws.Range("C31:C33").PasteSpecial(-4122)

# Row 31: Venice
$ws.Range("A31").Value = "Venice"
$ws.Range("B31").Value = "MS River at Venice"
$ws.Range("C31").Value = 0

# Row 32: BayouSale
$ws.Range("A32").Value = "BayouSale"
$ws.Range("B32").Value = "GIWW at Bayou Sale Ridge"
$ws.Range("C32").Value = 0

# Row 33: BayouBoeuf
$ws.Range("A33").Value = "BayouBoeuf"
$ws.Range("B33").Value = "Bayou Boeuf at Railroad Bridge"
$ws.Range("C33").Value = 0

# --- 3) Update the sheet's active selection to sit just below the ---
#        newly added data, matching the author's saved cursor
#        position.
[void]$ws.Range("A34").Select()
